$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 9794.24
$ws.Range("B8").Value = 9835.5499999999993
$ws.Range("C8").Value = 80.11
$ws.Range("D8").Value = 79.77
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -0.42
$ws.Range("G8").Value = 42609.488333333335
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("H8").Value = $false
